# Update LR-pair TPM-derived statistics with newly computed TPM values.
# Only numeric columns G,H,I,J,K,L,M,N,O,P,Q,R,S,T on rows 2-10 change;
# columns A-F are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 23.000594
$ws.Range("H2").Value = 69.00178200000001
$ws.Range("I2").Value = 0.1085495600721081
$ws.Range("J2").Value = 0.1085495600721081
$ws.Range("M2").Value = 32.63563666666667
$ws.Range("N2").Value = 97.90691000000001
$ws.Range("O2").Value = 0.9900101876891448
$ws.Range("P2").Value = 0.9900101876891446
$ws.Range("Q2").Value = 750.6390289015135
$ws.Range("R2").Value = 6755.751260113621
$ws.Range("S2").Value = 0.1074651703405619
$ws.Range("T2").Value = 0.1074651703405619

# Row 3
$ws.Range("G3").Value = 23.000594
$ws.Range("H3").Value = 69.00178200000001
$ws.Range("I3").Value = 0.1085495600721081
$ws.Range("J3").Value = 0.1085495600721081
$ws.Range("M3").Value = 0.05920466666666666
$ws.Range("O3").Value = 0.001795988347259859
$ws.Range("P3").Value = 0.001795988347259859
$ws.Range("Q3").Value = 1.361742500905333
$ws.Range("R3").Value = 12.255682508148
$ws.Range("S3").Value = 0.0001949537449896903
$ws.Range("T3").Value = 0.0001949537449896902

# Row 4
$ws.Range("G4").Value = 23.000594
$ws.Range("H4").Value = 69.00178200000001
$ws.Range("I4").Value = 0.1085495600721081
$ws.Range("J4").Value = 0.1085495600721081
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.270109
$ws.Range("N4").Value = 0.810327
$ws.Range("O4").Value = 0.008193823963595435
$ws.Range("P4").Value = 0.008193823963595434
$ws.Range("Q4").Value = 6.212667444746001
$ws.Range("R4").Value = 55.91400700271401
$ws.Range("S4").Value = 0.0008894359865565819
$ws.Range("T4").Value = 0.0008894359865565818

# Row 5
$ws.Range("G5").Value = 161.0956266666667
$ws.Range("H5").Value = 483.28688
$ws.Range("I5").Value = 0.7602785999442988
$ws.Range("J5").Value = 0.7602785999442987
$ws.Range("M5").Value = 32.63563666666667
$ws.Range("N5").Value = 97.90691000000001
$ws.Range("O5").Value = 0.9900101876891448
$ws.Range("P5").Value = 0.9900101876891446
$ws.Range("Q5").Value = 5257.458340482312
$ws.Range("R5").Value = 47317.12506434081
$ws.Range("S5").Value = 0.7526835594268956
$ws.Range("T5").Value = 0.7526835594268952

# Row 6
$ws.Range("G6").Value = 161.0956266666667
$ws.Range("H6").Value = 483.28688
$ws.Range("I6").Value = 0.7602785999442988
$ws.Range("J6").Value = 0.7602785999442987
$ws.Range("M6").Value = 0.05920466666666666
$ws.Range("O6").Value = 0.001795988347259859
$ws.Range("P6").Value = 0.001795988347259859
$ws.Range("Q6").Value = 9.537612878257777
$ws.Range("R6").Value = 85.83851590431999
$ws.Range("S6").Value = 0.001365451506171
$ws.Range("T6").Value = 0.001365451506171

# Row 7
$ws.Range("G7").Value = 161.0956266666667
$ws.Range("H7").Value = 483.28688
$ws.Range("I7").Value = 0.7602785999442988
$ws.Range("J7").Value = 0.7602785999442987
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.270109
$ws.Range("N7").Value = 0.810327
$ws.Range("O7").Value = 0.008193823963595435
$ws.Range("P7").Value = 0.008193823963595434
$ws.Range("Q7").Value = 43.51337862330666
$ws.Range("R7").Value = 391.62040760976
$ws.Range("S7").Value = 0.006229589011232383
$ws.Range("T7").Value = 0.00622958901123238

# Row 8
$ws.Range("G8").Value = 27.79403466666666
$ws.Range("H8").Value = 83.382104
$ws.Range("I8").Value = 0.131171839983593
$ws.Range("J8").Value = 0.131171839983593
$ws.Range("M8").Value = 32.63563666666667
$ws.Range("N8").Value = 97.90691000000001
$ws.Range("O8").Value = 0.9900101876891448
$ws.Range("P8").Value = 0.9900101876891446
$ws.Range("Q8").Value = 907.0760168820711
$ws.Range("R8").Value = 8163.684151938641
$ws.Range("S8").Value = 0.1298614579216874
$ws.Range("T8").Value = 0.1298614579216874

# Row 9
$ws.Range("G9").Value = 27.79403466666666
$ws.Range("H9").Value = 83.382104
$ws.Range("I9").Value = 0.131171839983593
$ws.Range("J9").Value = 0.131171839983593
$ws.Range("M9").Value = 0.05920466666666666
$ws.Range("O9").Value = 0.001795988347259859
$ws.Range("P9").Value = 0.001795988347259859
$ws.Range("Q9").Value = 1.645536557761778
$ws.Range("R9").Value = 14.809829019856
$ws.Range("S9").Value = 0.0002355830960991679
$ws.Range("T9").Value = 0.0002355830960991678

# Row 10
$ws.Range("G10").Value = 27.79403466666666
$ws.Range("H10").Value = 83.382104
$ws.Range("I10").Value = 0.131171839983593
$ws.Range("J10").Value = 0.131171839983593
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.270109
$ws.Range("N10").Value = 0.810327
$ws.Range("O10").Value = 0.008193823963595435
$ws.Range("P10").Value = 0.008193823963595434
$ws.Range("Q10").Value = 7.507418909778666
$ws.Range("R10").Value = 67.56677018800799
$ws.Range("S10").Value = 0.00107479896580647
$ws.Range("T10").Value = 0.00107479896580647
